$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

$ws.Range("D2").Value = '51.200.23'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.915.92'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue $ws.Range("D5") '363.86'
$ws.Range("E5").Value = '  +2.21%  '
Set-TextValue $ws.Range("D6") '104.69'
$ws.Range("E6").Value = '  -4.47%  '
Set-TextValue $ws.Range("D7") '0.541'
$ws.Range("E7").Value = '  -4.94%  '
$ws.Range("E8").Value = '  +0.05%  '
Set-TextValue $ws.Range("D9") '0.590'
$ws.Range("E9").Value = '  -6.30%  '
Set-TextValue $ws.Range("D10") '36.97'
$ws.Range("E10").Value = '  -4.88%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("E12").Value = '  -4.15%  '
Set-TextValue $ws.Range("D13") '18.46'
$ws.Range("E13").Value = '  -5.33%  '
$ws.Range("D14").Value = '3.371.70'
$ws.Range("E14").Value = '  +0.13%  '
Set-TextValue $ws.Range("D15") '7.35'
$ws.Range("D16").Value = '2.907.62'
$ws.Range("E16").Value = '  -0.29%  '
Set-TextValue $ws.Range("D17") '0.953'
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("D18").Value = '51.124.61'
$ws.Range("E18").Value = '  -1.58%  '
Set-TextValue $ws.Range("D19") '3.31'
$ws.Range("E19").Value = '  -3.89%  '
Set-TextValue $ws.Range("D20") '7.24'
$ws.Range("E20").Value = '  -3.98%  '
Set-TextValue $ws.Range("D21") '13.01'
$ws.Range("E21").Value = '  -6.55%  '
$ws.Range("D22").Value = '0.0₃0947'
$ws.Range("E22").Value = '  -3.23%  '
Set-TextValue $ws.Range("D23") '68.21'
$ws.Range("E23").Value = '  -3.22%  '
Set-TextValue $ws.Range("D24") '259.32'
$ws.Range("E24").Value = '  -3.35%  '
Set-TextValue $ws.Range("D25") '2.68'
$ws.Range("E25").Value = '  -4.44%  '
Set-TextValue $ws.Range("D26") '0.173'
$ws.Range("E26").Value = '  -5.17%  '
$ws.Range("E27").Value = '  +0.02%  '
Set-TextValue $ws.Range("D28") '25.91'
$ws.Range("E28").Value = '  -3.58%  '
Set-TextValue $ws.Range("D29") '7.17'
$ws.Range("E29").Value = '  -6.24%  '
$ws.Range("E30").Value = '  +0.40%  '
Set-TextValue $ws.Range("D31") '6.17'
$ws.Range("E31").Value = '  -0.19%  '
Set-TextValue $ws.Range("D32") '9.92'
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("E33").Value = '  -3.17%  '
Set-TextValue $ws.Range("D34") '34.96'
$ws.Range("E34").Value = '  -6.88%  '
Set-TextValue $ws.Range("D35") '50.71'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("E36").Value = '  +0.24%  '
Set-TextValue $ws.Range("D37") '0.0421'
$ws.Range("E37").Value = '  -4.79%  '
$ws.Range("E38").Value = '  +3.25%  '
Set-TextValue $ws.Range("D39") '3.14'
$ws.Range("E39").Value = '  -1.97%  '
Set-TextValue $ws.Range("D40") '16.95'
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("E41").Value = '  -6.51%  '
Set-TextValue $ws.Range("D42") '0.114'
$ws.Range("E42").Value = '  -4.65%  '
Set-TextValue $ws.Range("D43") '22.35'
$ws.Range("E43").Value = '  -2.22%  '
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("D46").Value = '2.062.94'
$ws.Range("E46").Value = '  -3.02%  '
Set-TextValue $ws.Range("D47") '3.19'
$ws.Range("E47").Value = '  -7.52%  '
Set-TextValue $ws.Range("D48") '2.27'
$ws.Range("E48").Value = '  -8.29%  '
$ws.Range("D49").Value = '3.194.28'
$ws.Range("E49").Value = '  -0.09%  '
Set-TextValue $ws.Range("D50") '0.236'
$ws.Range("E50").Value = '  -6.78%  '
Set-TextValue $ws.Range("D51") '0.0312'
$ws.Range("E51").Value = '  -7.53%  '
